# Various minor improvements in paper text, code, style, and refs, #4
#
# Adjust paragraph spacing on a handful of built-in / custom styles so
# headings and captions breathe a bit more.

$d = $word.ActiveDocument

# Heading 2: keep space-before (238 twips / 11.9pt), add space-after (40 twips / 2pt)
$style = $d.Styles("Heading 2")
$style.ParagraphFormat.SpaceBefore = 11.9
$style.ParagraphFormat.SpaceAfter = 2

# Heading 3: bump space-before to match Heading 2 (238 twips / 11.9pt),
# add space-after (40 twips / 2pt)
$style = $d.Styles("Heading 3")
$style.ParagraphFormat.SpaceBefore = 11.9
$style.ParagraphFormat.SpaceAfter = 2

# Image Caption: keep space-before (119 twips / 5.95pt), increase
# space-after from 119 to 181 twips (5.95pt -> 9.05pt)
$style = $d.Styles("Image Caption")
$style.ParagraphFormat.SpaceBefore = 5.95
$style.ParagraphFormat.SpaceAfter = 9.05

# References-title: keep space-before (181 twips / 9.05pt), add
# space-after (181 twips / 9.05pt)
$style = $d.Styles("References-title")
$style.ParagraphFormat.SpaceBefore = 9.05
$style.ParagraphFormat.SpaceAfter = 9.05
